$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 105, shifting existing rows 105+ down by one
$ws.Rows("105:105").Insert()

# Populate the newly inserted row 105 with the new stored procedure entry
$ws.Range("A105").Value = "stored procedure"
$ws.Range("B105").Value = "modify_product_group"
$ws.Range("C105").Value = "procedure for inserting or modifying product group"

# Leave the selection where data entry finished, matching the viewport
# the author ended up in after typing the new row.
[void]$ws.Range("D105").Select()
